# Actualización automática 2025-09-18 09:12:30
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M10").Value = 226.8
$ws1.Range("M20").Value = 1128.16
$ws1.Range("O20").Value = 261.27
$ws1.Range("P20").Value = 289.83
$ws1.Range("D28").Value = 457.92
$ws1.Range("D35").Value = "3 de 33"
$ws1.Range("M35").Value = "7 de 33"
$ws1.Range("O35").Value = "1 de 33"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F10").Value = 226.8
$ws2.Range("F20").Value = 1825.21
$ws2.Range("F28").Value = 457.92
$ws2.Range("F35").Value = 15221.26

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D3").Value = 2218.75
$ws3.Range("E3").Value = 6615.82354940916
$ws3.Range("F3").Value = 0.2511439842105775

$ws3.Range("D8").Value = 289.83
$ws3.Range("E8").Value = 190.387743214072
$ws3.Range("F8").Value = 0.6035387157088014

$ws3.Range("D12").Value = 9338.65
$ws3.Range("E12").Value = 13095.1053751766
$ws3.Range("F12").Value = 0.4162767153257543

$ws3.Range("D14").Value = 261.27
$ws3.Range("E14").Value = 1163.6962010375
$ws3.Range("F14").Value = 0.1833517172616253

$ws3.Range("D15").Value = 15478.49
$ws3.Range("E15").Value = 23264.52881339593
$ws3.Range("F15").Value = 0.3995168800488025
